$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range (A1:G4) so stale columns F/G and any
# leftover cells from the old, shorter table are removed before writing
# the new layout.
$ws.Range("A1:G11").ClearContents()

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (Índice, Distancia, max, min, Tempo)
$data = @(
    @(0, 7412.066666666667, 8011, 6888, 0.1389382123947143),
    @(1, 7914.2,            8561, 7472, 0.1821735779444377),
    @(2, 8134.533333333334, 8709, 7353, 0.199787974357605),
    @(3, 8661.866666666667, 9485, 7842, 0.1983153740564982),
    @(4, 7755.366666666667, 8157, 7334, 0.1859331528345744),
    @(5, 7961.333333333333, 8621, 7428, 0.161428181330363),
    @(6, 7914.7,            8429, 7129, 0.1499676465988159),
    @(7, 7682.966666666666, 8354, 6682, 0.1341055949529012),
    @(8, 7146.166666666667, 7599, 6387, 0.1559197028477987),
    @(9, 7523.2,            7979, 6389, 0.1735461393992106)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
